# Add an "exclude prioritized" Priority column to the Screen Print Designs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Screen Print Designs")

# Insert a new column before column F (Subcategory1 col), shifting everything
# from F onward one column to the right.
$ws.Columns("F:F").Insert()

# New column header + the one data point the author filled in.
$ws.Range("F1").Value = "Priority"
$ws.Range("F13").Value = 1

$ws.Range("E27").Select() | Out-Null
